$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts the string into a
# numeric value (losing formatting like trailing zeros).
$textCells = @("D5", "D6", "D9", "D10", "D12", "D14", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D33", "D34", "D36", "D37", "D40", "D41", "D42", "D46", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.372.44"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.805.06"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "666.16"
$ws.Range("E5").Value = "  +7.00%  "
$ws.Range("D6").Value = "168.34"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "3.803.24"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "7.01"
$ws.Range("E12").Value = "  +5.30%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "35.77"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "4.451.99"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "3.807.79"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "70.360.42"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "17.70"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "7.15"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "10.78"
$ws.Range("E21").Value = "  +12.24%  "
$ws.Range("D22").Value = "474.06"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").Value = "0.712"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "82.72"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "0.0000143"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "10.32"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "3.958.17"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +7.19%  "
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("D33").Value = "7.39"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "29.63"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +11.83%  "
$ws.Range("D36").Value = "9.14"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.764.12"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").Value = "0.970"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +9.25%  "
$ws.Range("D46").Value = "45.79"
$ws.Range("E46").Value = "  +6.00%  "
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D48").Value = "157.30"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "0.300"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "1.42"
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  +0.97%  "
